$wb = $excel.ActiveWorkbook

$hotel  = $wb.Worksheets.Item("hotel_info")
$review = $wb.Worksheets.Item("review_info")

# 1. Insert a new "State" column into hotel_info between Hotel_Name and City,
#    and populate the header + the value for the existing data row.
$hotel.Range("C1").EntireColumn.Insert()
$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# 2. Reorder the worksheet tabs so review_info comes before hotel_info.
$review.Move($hotel)
